$wb = $excel.ActiveWorkbook
$ws = $wb.ActiveSheet

# Updated values for existing rows 218-222 (columns B..U), reflecting the
# monthly data revision that came with the new period below.
$updates = @{
    218 = @(210862,2889,209,2680,5410,2,5408,0,0,32534,7464,25070,70941,11051,59889,97911,1445,96465,1177,1177)
    219 = @(209836,2840,176,2664,5294,2,5293,0,0,32171,7390,24781,71830,11511,60319,96525,2090,94435,1176,1176)
    220 = @(209510,2829,174,2655,5255,2,5253,0,0,33354,8504,24850,72606,12382,60224,94308,2060,92249,1158,1158)
    221 = @(213397,2729,100,2629,5115,2,5113,0,0,33299,8473,24826,73639,12335,61305,97441,1690,95751,1173,1173)
    222 = @(214870,2684,73,2612,5075,2,5073,0,0,33814,9456,24357,73114,12489,60626,98989,1445,97544,1193,1193)
}

foreach ($r in $updates.Keys) {
    $vals = $updates[$r]
    for ($i = 0; $i -lt $vals.Length; $i++) {
        # Column B is index 2 (A=1), so offset by +2
        $ws.Cells.Item($r, $i + 2).Value = $vals[$i]
    }
}

# New row 223: new monthly period "01-06-2021".
# Build the label through a formula + paste-as-values round trip so Excel
# stores it as a genuine shared-string cell instead of silently recognising
# it as a date literal (which would reformat the cell and grow styles.xml).
$newRow = 223
$labelCell = $ws.Range("A$newRow")
$labelCell.Formula = '="01-06-2021"'
$labelCell.Copy() | Out-Null
$labelCell.PasteSpecial(-4163) | Out-Null  # xlPasteValues
$excel.CutCopyMode = $false

$newRowValues = @(213283,2664,73,2591,4955,2,4954,0,0,32910,8785,24125,71901,12272,59629,99687,1186,98501,1165,1165)
for ($i = 0; $i -lt $newRowValues.Length; $i++) {
    $ws.Cells.Item($newRow, $i + 2).Value = $newRowValues[$i]
}
